$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("_input")
Write-Host "Sheet name: " $ws.Name
